$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$rng = $ws.Range("A1:E31")
$rng.Borders.Item(1).Weight = 4
$rng.Borders.Item(1).Color = 255
